$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '70.946.29'
$ws.Range('E2').Value = '  +0.52%  '

$ws.Range('D3').Value = '3.543.99'
$ws.Range('E3').Value = '  -0.52%  '

$ws.Range('E4').Value = '  +0.00%  '

$r = $ws.Range('D5')
$r.NumberFormat = '@'
$r.Value = '618.09'
$r.Style = 'Normal'
$ws.Range('E5').Value = '  +1.01%  '

$r = $ws.Range('D6')
$r.NumberFormat = '@'
$r.Value = '175.12'
$r.Style = 'Normal'
$ws.Range('E6').Value = '  +1.15%  '

$ws.Range('D7').Value = '3.539.14'
$ws.Range('E7').Value = '  -0.40%  '

$ws.Range('E8').Value = '  -1.01%  '

$ws.Range('E9').Value = '  -0.04%  '

$r = $ws.Range('D10')
$r.NumberFormat = '@'
$r.Value = '0.200'
$r.Style = 'Normal'
$ws.Range('E10').Value = '  +1.64%  '

$r = $ws.Range('D11')
$r.NumberFormat = '@'
$r.Value = '7.25'
$r.Style = 'Normal'
$ws.Range('E11').Value = '  -2.98%  '

$ws.Range('E12').Value = '  +0.52%  '

$r = $ws.Range('D13')
$r.NumberFormat = '@'
$r.Value = '46.79'
$r.Style = 'Normal'
$ws.Range('E13').Value = '  +0.27%  '

$ws.Range('E14').Value = '  +0.07%  '

$ws.Range('D15').Value = '4.116.74'
$ws.Range('E15').Value = '  -0.53%  '

$r = $ws.Range('D16')
$r.NumberFormat = '@'
$r.Value = '8.46'
$r.Style = 'Normal'
$ws.Range('E16').Value = '  +0.97%  '

$r = $ws.Range('D17')
$r.NumberFormat = '@'
$r.Value = '612.81'
$r.Style = 'Normal'
$ws.Range('E17').Value = '  -0.50%  '

$ws.Range('D18').Value = '3.559.76'
$ws.Range('E18').Value = '  +0.03%  '

$ws.Range('D19').Value = '71.048.39'
$ws.Range('E19').Value = '  +0.54%  '

$ws.Range('E20').Value = '  +1.14%  '

$r = $ws.Range('D21')
$r.NumberFormat = '@'
$r.Value = '17.78'
$r.Style = 'Normal'
$ws.Range('E21').Value = '  +2.23%  '

$r = $ws.Range('D22')
$r.NumberFormat = '@'
$r.Value = '0.890'
$r.Style = 'Normal'
$ws.Range('E22').Value = '  +0.71%  '

$r = $ws.Range('D23')
$r.NumberFormat = '@'
$r.Value = '9.09'
$r.Style = 'Normal'
$ws.Range('E23').Value = '  -2.89%  '

$r = $ws.Range('D24')
$r.NumberFormat = '@'
$r.Value = '15.76'
$r.Style = 'Normal'
$ws.Range('E24').Value = '  -1.87%  '

$r = $ws.Range('D25')
$r.NumberFormat = '@'
$r.Value = '98.61'
$r.Style = 'Normal'
$ws.Range('E25').Value = '  +1.52%  '

$ws.Range('E26').Value = '  -0.93%  '

$ws.Range('E27').Value = '  +0.07%  '

$ws.Range('E28').Value = '  -0.91%  '

$r = $ws.Range('D29')
$r.NumberFormat = '@'
$r.Value = '33.93'
$r.Style = 'Normal'
$ws.Range('E29').Value = '  +1.42%  '

$r = $ws.Range('D30')
$r.NumberFormat = '@'
$r.Value = '9.17'
$r.Style = 'Normal'
$ws.Range('E30').Value = '  +1.11%  '

$r = $ws.Range('D31')
$r.NumberFormat = '@'
$r.Value = '3.05'
$r.Style = 'Normal'
$ws.Range('E31').Value = '  +0.07%  '

$ws.Range('E32').Value = '  -3.95%  '

$ws.Range('E33').Value = '  +0.14%  '

$ws.Range('E34').Value = '  -1.18%  '

$r = $ws.Range('D35')
$r.NumberFormat = '@'
$r.Value = '632.51'
$r.Style = 'Normal'
$ws.Range('E35').Value = '  +9.76%  '

$ws.Range('E36').Value = '  -0.91%  '

$r = $ws.Range('D37')
$r.NumberFormat = '@'
$r.Value = '10.87'
$r.Style = 'Normal'
$ws.Range('E37').Value = '  +0.47%  '

$ws.Range('E38').Value = '  -4.78%  '

$r = $ws.Range('D39')
$r.NumberFormat = '@'
$r.Value = '0.0479'
$r.Style = 'Normal'
$ws.Range('E39').Value = '  -0.05%  '

$r = $ws.Range('D40')
$r.NumberFormat = '@'
$r.Value = '57.04'
$r.Style = 'Normal'
$ws.Range('E40').Value = '  -0.48%  '

$ws.Range('E41').Value = '  +0.09%  '

$ws.Range('E42').Value = '  +2.25%  '

$ws.Range('D43').Value = '0.0₃0745'
$ws.Range('E43').Value = '  +5.47%  '

$ws.Range('D44').Value = '3.374.80'
$ws.Range('E44').Value = '  -0.49%  '

$r = $ws.Range('D45')
$r.NumberFormat = '@'
$r.Value = '3.00'
$r.Style = 'Normal'
$ws.Range('E45').Value = '  +0.20%  '

$ws.Range('E46').Value = '  -1.98%  '

$r = $ws.Range('D47')
$r.NumberFormat = '@'
$r.Value = '32.25'
$r.Style = 'Normal'
$ws.Range('E47').Value = '  -2.82%  '

$ws.Range('E48').Value = '  -1.50%  '

$ws.Range('E49').Value = '  +0.51%  '

$r = $ws.Range('D50')
$r.NumberFormat = '@'
$r.Value = '133.24'
$r.Style = 'Normal'
$ws.Range('E50').Value = '  -0.32%  '

$ws.Range('E51').Value = '  +6.40%  '

